$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data

# Row 2
$ws.Range('D2').Value = '26.155.05'
$ws.Range('E2').Value = '  +1.42%  '

# Row 3
$ws.Range('D3').Value = '1.612.81'
$ws.Range('E3').Value = '  +0.92%  '

# Row 4
$ws.Range('E4').Value = '  -0.50%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.14%  '

# Row 6
$ws.Range('E6').Value = '  -0.48%  '

# Row 7
$ws.Range('E7').Value = '  +1.62%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.81%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.34%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.39'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.43%  '

# Row 11
$ws.Range('E11').Value = '  +1.42%  '

# Row 12
$ws.Range('D12').Value = '1.837.10'
$ws.Range('E12').Value = '  +0.91%  '

# Row 13
$ws.Range('D13').Value = '1.611.56'
$ws.Range('E13').Value = '  +0.73%  '

# Row 14
$ws.Range('E14').Value = '  +0.61%  '

# Row 15
$ws.Range('E15').Value = '  +0.84%  '

# Row 16
$ws.Range('D16').Value = '26.155.96'
$ws.Range('E16').Value = '  +1.43%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.25%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  +2.25%  '

# Row 19
$ws.Range('E19').Value = '  -0.44%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '199.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.60%  '

# Row 21
$ws.Range('E21').Value = '  +2.49%  '

# Row 22
$ws.Range('E22').Value = '  +2.51%  '

# Row 23
$ws.Range('E23').Value = '  +1.83%  '

# Row 24
$ws.Range('E24').Value = '  +4.19%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.50%  '

# Row 27
$ws.Range('E27').Value = '  -0.51%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.61%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.18%  '

# Row 30
$ws.Range('E30').Value = '  -1.38%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0476'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.07%  '

# Row 32
$ws.Range('E32').Value = '  +2.50%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.24%  '

# Row 34
$ws.Range('E34').Value = '  +4.02%  '

# Row 35
$ws.Range('E35').Value = '  -1.89%  '

# Row 36
$ws.Range('D36').Value = '1.108.88'
$ws.Range('E36').Value = '  +1.48%  '

# Row 37
$ws.Range('E37').Value = '  +1.73%  '

# Row 38
$ws.Range('E38').Value = '  -0.49%  '

# Row 39
$ws.Range('E39').Value = '  +3.49%  '

# Row 40
$ws.Range('E40').Value = '  -1.01%  '

# Row 41
$ws.Range('E41').Value = '  +0.16%  '

# Row 42
$ws.Range('E42').Value = '  +8.01%  '

# Row 43
$ws.Range('D43').Value = '1.748.65'
$ws.Range('E43').Value = '  +0.90%  '

# Row 44
$ws.Range('E44').Value = '  +1.71%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.48%  '

# Row 46
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.40%  '

# Row 47
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0105'
$ws.Range('E47').Value = '  +8.03%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.84%  '

# Row 49
$ws.Range('E49').Value = '  +0.02%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.408'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.71%  '

# Row 51
$ws.Range('E51').Value = '  -0.53%  '
